$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column C for rows 2-5 from 2023-11-13 to 2023-11-14
$ws.Range("C2:C5").Value = 45244
